$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

$ws.Range("D2").Value = "36.623.19"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "1.963.03"
$ws.Range("E3").Value = "  +0.73%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue "D5" "244.38"
$ws.Range("E5").Value = "  +0.30%  "

Set-TextValue "D6" "0.619"
$ws.Range("E6").Value = "  +0.26%  "

Set-TextValue "D7" "58.99"
$ws.Range("E7").Value = "  +1.18%  "

$ws.Range("E8").Value = "  +0.00%  "

Set-TextValue "D9" "0.373"
$ws.Range("E9").Value = "  +1.68%  "

Set-TextValue "D10" "0.0815"
$ws.Range("E10").Value = "  -3.00%  "

$ws.Range("E11").Value = "  -0.57%  "

Set-TextValue "D12" "22.18"
$ws.Range("E12").Value = "  +2.41%  "

$ws.Range("D13").Value = "2.252.76"
$ws.Range("E13").Value = "  +0.79%  "

Set-TextValue "D14" "0.828"
$ws.Range("E14").Value = "  +0.26%  "

Set-TextValue "D15" "13.70"
$ws.Range("E15").Value = "  +0.68%  "

Set-TextValue "D16" "5.26"
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("D17").Value = "1.982.85"
$ws.Range("E17").Value = "  +2.48%  "

$ws.Range("D18").Value = "36.525.52"
$ws.Range("E18").Value = "  +0.27%  "

Set-TextValue "D19" "69.95"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").Value = "0.0₃0858"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "229.08"
$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D22" "5.07"
$ws.Range("E22").Value = "  -0.10%  "

Set-TextValue "D23" "0.999"
$ws.Range("E23").Value = "  -0.17%  "

Set-TextValue "D24" "2.43"
$ws.Range("E24").Value = "  -0.51%  "

Set-TextValue "D25" "2.36"
$ws.Range("E25").Value = "  +2.84%  "

Set-TextValue "D26" "0.141"
$ws.Range("E26").Value = "  +13.21%  "

Set-TextValue "D27" "9.24"
$ws.Range("E27").Value = "  -0.08%  "

Set-TextValue "D28" "160.16"
$ws.Range("E28").Value = "  -1.41%  "

Set-TextValue "D29" "19.37"
$ws.Range("E29").Value = "  -0.30%  "

Set-TextValue "D30" "0.120"
$ws.Range("E30").Value = "  +1.53%  "

Set-TextValue "D31" "1.15"
$ws.Range("E31").Value = "  -0.42%  "

Set-TextValue "D32" "4.70"
$ws.Range("E32").Value = "  +0.42%  "

Set-TextValue "D33" "0.0617"
$ws.Range("E33").Value = "  -1.88%  "

Set-TextValue "D34" "4.28"
$ws.Range("E34").Value = "  -0.30%  "

Set-TextValue "D35" "2.28"
$ws.Range("E35").Value = "  +6.28%  "

$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D37" "3.39"
$ws.Range("E37").Value = "  +10.94%  "

$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D38" "6.01"
$ws.Range("E38").Value = "  -3.57%  "

$ws.Range("E39").Value = "  -0.42%  "

Set-TextValue "D40" "0.0985"
$ws.Range("E40").Value = "  -0.28%  "

$ws.Range("E41").Value = "  +0.44%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D42" "0.0212"
$ws.Range("E42").Value = "  +0.92%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D43" "1.17"
$ws.Range("E43").Value = "  -0.85%  "

Set-TextValue "D44" "16.18"

$ws.Range("D45").Value = "1.359.20"
$ws.Range("E45").Value = "  +0.59%  "

$ws.Range("E46").Value = "  +0.45%  "

Set-TextValue "D47" "87.75"
$ws.Range("E47").Value = "  -0.18%  "

Set-TextValue "D48" "7.14"
$ws.Range("E48").Value = "  -0.62%  "

Set-TextValue "D49" "2.83"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("D50").Value = "2.144.14"
$ws.Range("E50").Value = "  +0.85%  "

Set-TextValue "D51" "43.73"
$ws.Range("E51").Value = "  -4.16%  "
